$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 3.5
$ws.Range("I3").Value = 2.4
$ws.Range("J3").Value = 1.14
$ws.Range("K3").Value = 5.5
$ws.Range("L3").Value = 1.57
$ws.Range("M3").Value = 2.25
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("T3").Value = 7.5
$ws.Range("U3").Value = 15
$ws.Range("X3").Value = 34
$ws.Range("AF3").Value = 10
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 23

# Row 4
$ws.Range("G4").Value = 2.25
$ws.Range("I4").Value = 3
$ws.Range("K4").Value = 10
$ws.Range("U4").Value = 10
$ws.Range("X4").Value = 19
$ws.Range("AE4").Value = 9
$ws.Range("AF4").Value = 15
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 26

# Row 5
$ws.Range("J5").Value = 1.02
$ws.Range("K5").Value = 21

# Row 6
$ws.Range("G6").Value = 1.71

# Row 8
$ws.Range("G8").Value = 1.31
$ws.Range("H8").Value = 4.65
$ws.Range("I8").Value = 9.25
$ws.Range("L8").Value = 1.26
$ws.Range("M8").Value = 3.2
$ws.Range("N8").Value = 1.75
$ws.Range("O8").Value = 1.85
$ws.Range("R8").Value = 2.22
$ws.Range("S8").Value = 1.52
$ws.Range("T8").Value = 5.8
$ws.Range("U8").Value = 5.5
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 7.4
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 37
$ws.Range("Z8").Value = 10.5
$ws.Range("AA8").Value = 9.5
$ws.Range("AB8").Value = 27
$ws.Range("AC8").Value = 175
$ws.Range("AE8").Value = 21
$ws.Range("AF8").Value = 65
$ws.Range("AI8").Value = 150
$ws.Range("AJ8").Value = 120

# Row 12
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 2.35
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 1.2
$ws.Range("K12").Value = 4.05
$ws.Range("L12").Value = 1.75
$ws.Range("M12").Value = 1.98
$ws.Range("N12").Value = 3.15
$ws.Range("O12").Value = 1.31
$ws.Range("P12").Value = 1.75
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 2.27
$ws.Range("S12").Value = 1.57
$ws.Range("T12").Value = 5.9
$ws.Range("U12").Value = 13.5
$ws.Range("V12").Value = 11.75
$ws.Range("W12").Value = 45
$ws.Range("X12").Value = 37
$ws.Range("Y12").Value = 65
$ws.Range("Z12").Value = 4.05
$ws.Range("AA12").Value = 5
$ws.Range("AB12").Value = 19
$ws.Range("AC12").Value = 150
$ws.Range("AE12").Value = 6
$ws.Range("AF12").Value = 14
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 45
$ws.Range("AJ12").Value = 65

# Row 13
$ws.Range("G13").Value = 2.27
$ws.Range("I13").Value = 3.6
$ws.Range("K13").Value = 4.5
$ws.Range("L13").Value = 1.72
$ws.Range("N13").Value = 3.1
$ws.Range("O13").Value = 1.32
$ws.Range("T13").Value = 4.85
$ws.Range("U13").Value = 8.75
$ws.Range("W13").Value = 23
$ws.Range("Z13").Value = 4.5
$ws.Range("AA13").Value = 5.9
$ws.Range("AB13").Value = 24
$ws.Range("AE13").Value = 6.8
$ws.Range("AF13").Value = 16.5
$ws.Range("AI13").Value = 50

# Row 16
$ws.Range("G16").Value = 2.1
$ws.Range("I16").Value = 3.5
$ws.Range("J16").Value = 1.07
$ws.Range("K16").Value = 9
$ws.Range("L16").Value = 1.36
$ws.Range("M16").Value = 3
$ws.Range("U16").Value = 9.5
$ws.Range("W16").Value = 19
$ws.Range("AH16").Value = 41

# Row 19
$ws.Range("J19").Value = 1.14
$ws.Range("K19").Value = 5.5
$ws.Range("Z19").Value = 5.5
$ws.Range("AB19").Value = 21
$ws.Range("AJ19").Value = 51

# Row 20
$ws.Range("J20").Value = 1.06
$ws.Range("K20").Value = 10
$ws.Range("N20").Value = 2.08
$ws.Range("O20").Value = 1.73

# Row 22
$ws.Range("AI22").Value = 23

# Row 23
$ws.Range("I23").Value = 3.25
$ws.Range("L23").Value = 1.44
$ws.Range("M23").Value = 2.63

# Row 28
$ws.Range("L28").Value = 1.5
$ws.Range("M28").Value = 2.5

# Row 33
$ws.Range("V33").Value = 11.25
$ws.Range("AG33").Value = 9

# Row 36
$ws.Range("G36").Value = 2.9
$ws.Range("H36").Value = 3.25
$ws.Range("I36").Value = 2.5
$ws.Range("J36").Value = 1.07
$ws.Range("K36").Value = 9
$ws.Range("R36").Value = 1.8
$ws.Range("S36").Value = 1.95
$ws.Range("Z36").Value = 9
$ws.Range("AA36").Value = 6
$ws.Range("AD36").Value = 251
$ws.Range("AE36").Value = 8
$ws.Range("AF36").Value = 12
$ws.Range("AG36").Value = 10
$ws.Range("AJ36").Value = 34

# Row 42
$ws.Range("K42").Value = 13
